# Se optimizó el método crearMesas del formulario InicializarMesas
#
# Extiende la tabla de "Enero" en la hoja Horas2018 con tres filas más
# (12, 14 y 15 de enero) y actualiza los totales/resumen para que cubran
# el nuevo rango de datos (hasta la fila 26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas2018")

# --- Nuevas filas de datos --------------------------------------------------
# (se escriben los valores/fórmulas ANTES de copiar el formato para que el
# motor de cálculo registre correctamente las nuevas dependencias)
$ws.Range("B23").Value = "Enero"
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 3

$ws.Range("B25").Value = "Enero"
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 5

$ws.Range("B26").Value = "Enero"
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 6

# --- Resumen: las fórmulas ahora cubren hasta la fila 26 --------------------
$ws.Range("H3").Formula = "=COUNT(C4:C26)"
$ws.Range("H4").Formula = "=SUM(D4:D26)"

# --- Formato: igual que las filas existentes de "Enero" (p.ej. fila 22) ----
$source = $ws.Range("B22:D22")
$source.Copy()
$ws.Range("B23:D23").PasteSpecial(-4122)
$ws.Range("B25:D25").PasteSpecial(-4122)
$ws.Range("B26:D26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selección final, igual que la dejó el autor tras editar ---------------
$ws.Range("D26").Select()
